$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.3854739455626746
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = -0.03962185825936195
$ws.Range("G2").Value = 0.01805676298785952

# Row 3
$ws.Range("B3").Value = -0.05480118223795406
$ws.Range("C3").Value = 1.288231335644898
$ws.Range("D3").Value = 0.006641498110503894
$ws.Range("E3").Value = 0.004705187883590671
$ws.Range("F3").Value = -0.1075615368144068
$ws.Range("G3").Value = 0.1589222979733107

# Row 4
$ws.Range("B4").Value = 0.05483535554023129
$ws.Range("C4").Value = -5.234930118604909
$ws.Range("D4").Value = -0.3655849929097255
$ws.Range("E4").Value = -0.03295839414375842
$ws.Range("F4").Value = 0.7021308375697177
$ws.Range("G4").Value = 0.09584828949574475

# Row 5
$ws.Range("B5").Value = 0.7346593316326704
$ws.Range("C5").Value = 3.111218448870295
$ws.Range("D5").Value = 0.4556031319157312
$ws.Range("E5").Value = 0.0504320718900517
$ws.Range("F5").Value = 1.727024331227208
$ws.Range("G5").Value = 0.03731861978146789
